$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update computed Cost/UnitCost for rows 2 and 3 ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 364.2001544999999
$schedule.Range("F2").Value = 8.029103935185184
$schedule.Range("E3").Value = 459.2079375
$schedule.Range("F3").Value = 30.37089533730159

# --- Sheet "Detailed": refresh the DateTime/Price/Type/Date/Pump_Status table ---
$detailed = $wb.Worksheets.Item("Detailed")

# Remember the formats used by the existing data columns so that any newly
# created rows (48 and 49) pick up the same number formatting as the rest
# of the table.
$dateTimeFormat = $detailed.Cells.Item(2, 1).NumberFormat
$dateFormat = $detailed.Cells.Item(2, 4).NumberFormat

$rows = @(
    @(2, 46042, 58.98372, "historical", 46042, "OFF"),
    @(3, 46042.02083333334, 57.06003, "historical", 46042, "OFF"),
    @(4, 46042.04166666666, 57.06003, "historical", 46042, "OFF"),
    @(5, 46042.0625, 57.06003, "historical", 46042, "OFF"),
    @(6, 46042.08333333334, 57.06003, "historical", 46042, "OFF"),
    @(7, 46042.10416666666, 63.16199, "historical", 46042, "OFF"),
    @(8, 46042.125, 63.19773, "historical", 46042, "OFF"),
    @(9, 46042.14583333334, 63.98176, "forecast", 46042, "OFF"),
    @(10, 46042.16666666666, 64.16161, "forecast", 46042, "ON"),
    @(11, 46042.1875, 63.65567, "forecast", 46042, "ON"),
    @(12, 46042.20833333334, 70.53139, "forecast", 46042, "ON"),
    @(13, 46042.22916666666, 65.09452, "forecast", 46042, "ON"),
    @(14, 46042.25, 75.59721, "forecast", 46042, "ON"),
    @(15, 46042.27083333334, 66.60306, "forecast", 46042, "ON"),
    @(16, 46042.29166666666, 57.06003, "forecast", 46042, "ON"),
    @(17, 46042.3125, 8.59503, "forecast", 46042, "ON"),
    @(18, 46042.33333333334, 0.679, "forecast", 46042, "ON"),
    @(19, 46042.35416666666, 0.00983, "forecast", 46042, "ON"),
    @(20, 46042.375, -1.1032, "forecast", 46042, "ON"),
    @(21, 46042.39583333334, -5.58973, "forecast", 46042, "ON"),
    @(22, 46042.41666666666, -6.78005, "forecast", 46042, "ON"),
    @(23, 46042.4375, -6, "forecast", 46042, "ON"),
    @(24, 46042.45833333334, -8.244210000000001, "forecast", 46042, "ON"),
    @(25, 46042.47916666666, -7.65069, "forecast", 46042, "ON"),
    @(26, 46042.5, -6.91063, "forecast", 46042, "ON"),
    @(27, 46042.52083333334, -9.139889999999999, "forecast", 46042, "ON"),
    @(28, 46042.54166666666, -9.99, "forecast", 46042, "ON"),
    @(29, 46042.5625, -7.95296, "forecast", 46042, "ON"),
    @(30, 46042.58333333334, -7.78933, "forecast", 46042, "ON"),
    @(31, 46042.60416666666, -7.99712, "forecast", 46042, "ON"),
    @(32, 46042.625, -6.77628, "forecast", 46042, "ON"),
    @(33, 46042.64583333334, -6.52464, "forecast", 46042, "ON"),
    @(34, 46042.66666666666, -0.93103, "forecast", 46042, "ON"),
    @(35, 46042.6875, -5.51, "forecast", 46042, "OFF"),
    @(36, 46042.70833333334, -5.01, "forecast", 46042, "OFF"),
    @(37, 46042.72916666666, 9.60627, "forecast", 46042, "OFF"),
    @(38, 46042.75, 13.92433, "forecast", 46042, "OFF"),
    @(39, 46042.77083333334, 46.71265, "forecast", 46042, "OFF"),
    @(40, 46042.79166666666, 57.3, "forecast", 46042, "OFF"),
    @(41, 46042.8125, 57.3, "forecast", 46042, "OFF"),
    @(42, 46042.83333333334, 60.98539, "forecast", 46042, "ON"),
    @(43, 46042.85416666666, 58.51073, "forecast", 46042, "ON"),
    @(44, 46042.875, 57.42395, "forecast", 46042, "ON"),
    @(45, 46042.89583333334, 57.77587, "forecast", 46042, "ON"),
    @(46, 46042.91666666666, 57.3, "forecast", 46042, "ON"),
    @(47, 46042.9375, 58.63896, "forecast", 46042, "ON"),
    @(48, 46042.95833333334, 61.10068, "forecast", 46042, "ON"),
    @(49, 46042.97916666666, 59.24692, "forecast", 46042, "ON")
)

foreach ($row in $rows) {
    $r = $row[0]
    $detailed.Cells.Item($r, 1).Value = $row[1]
    $detailed.Cells.Item($r, 1).NumberFormat = $dateTimeFormat
    $detailed.Cells.Item($r, 2).Value = $row[2]
    $detailed.Cells.Item($r, 3).Value = $row[3]
    $detailed.Cells.Item($r, 4).Value = $row[4]
    $detailed.Cells.Item($r, 4).NumberFormat = $dateFormat
    $detailed.Cells.Item($r, 5).Value = $row[5]
}

